# devicelist.xlsx edit
#  1. Add a new "Browsers" worksheet with Sauce Labs Windows/Chrome rows
#  2. Bump testobject_appium_version (col H) from 1.9.1 -> 1.17.0 on RealDevices
#  3. Clear the stray platformVersion "11" on the ios_phone_only row (F7)
#  4. Leave the new sheet active/selected, matching the author's final UI state

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- add the Browsers sheet, right after RealDevices --------------------
$browsers = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$browsers.Name = "Browsers"

$browsers.Range("A1").Value = "name"
$browsers.Range("B1").Value = "platformName"
$browsers.Range("C1").Value = "platform"
$browsers.Range("D1").Value = "browser"
$browsers.Range("E1").Value = "version"
$browsers.Range("F1").Value = "platformVersion"
$browsers.Range("G1").Value = "chromeArguments"
$browsers.Range("H1").Value = "firefoxArguments"
$browsers.Range("I1").Value = "screenResolution"

# --- RealDevices: appium client bump + blank out the leftover "11" ------
$ws.Range("H2").Value = "'1.17.0"
$ws.Range("H3").Value = "'1.17.0"
$ws.Range("H4").Value = "'1.17.0"
$ws.Range("H5").Value = "'1.17.0"
$ws.Range("H6").Value = "'1.17.0"
$ws.Range("H7").Value = "'1.17.0"
$ws.Range("F7").Value = "'"
$ws.Range("A7").Select()

# --- Browsers row 2 (typed out in full, then duplicated downward) -------
$browsers.Range("B2").Value = "Windows 10"
$browsers.Range("C2").Value = "Windows 10"
$browsers.Range("D2").Value = "chrome"
$browsers.Range("E2").Value = "'80"
$browsers.Range("F2").Value = "'80"
$browsers.Range("G2").Value = "'"
$browsers.Range("H2").Value = "'"
$browsers.Range("I2").Value = "1280x960"
$browsers.Range("A2").Value = "SL_WIN10_CHROME_80"

$browsers.Range("A2:I2").Copy() | Out-Null
$browsers.Range("A3:I3").PasteSpecial(-4104) | Out-Null
$browsers.Range("A4:I4").PasteSpecial(-4104) | Out-Null
$browsers.Range("A5:I5").PasteSpecial(-4104) | Out-Null

# blank "arguments" cells don't survive PasteSpecial, restore them
$browsers.Range("G3").Value = "'"
$browsers.Range("H3").Value = "'"
$browsers.Range("G4").Value = "'"
$browsers.Range("H4").Value = "'"
$browsers.Range("G5").Value = "'"
$browsers.Range("H5").Value = "'"

# unique per-row name + chrome version for the other three rows
$browsers.Range("A3").Value = "SL_WIN10_CHROME_79"
$browsers.Range("A4").Value = "SL_WIN10_CHROME_78"
$browsers.Range("A5").Value = "SL_WIN10_CHROME_77"

$browsers.Range("E3").Value = "'79"
$browsers.Range("F3").Value = "'79"
$browsers.Range("E4").Value = "'78"
$browsers.Range("F4").Value = "'78"
$browsers.Range("E5").Value = "'77"
$browsers.Range("F5").Value = "'77"

# --- formatting: column widths -------------------------------------------
$browsers.Columns.Item(1).ColumnWidth = 24.16
$browsers.Range("B1:I5").ColumnWidth = 19
$browsers.Columns.Item(10).ColumnWidth = 10.66

# --- final UI state: Browsers tab active, B5 selected --------------------
$browsers.Range("B5").Select()
$browsers.Activate()
